$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": bump header timestamps/counter and append the new scrape
# row (row 74) for the 06:52:38 run.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:52:38"
$ws1.Range("A3").Value = "Total filas: 69"

$ws1.Cells.Item(74, 1).Value = "06:52:38"
$ws1.Cells.Item(74, 2).Value = "08:42"
$ws1.Cells.Item(74, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(74, 4).Value = 110
$ws1.Cells.Item(74, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": only the "last updated" timestamp moves forward.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:52:38"

# ---------------------------------------------------------------------------
# Sheet "6203-6173": bump header timestamp/counter, insert a fresh scrape
# row before the old row 15 (shifting the two existing rows down), and
# append a brand-new row at the end.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:52:38"
$ws3.Range("A3").Value = "Total filas: 13"

$ws3.Rows(15).Insert()
$ws3.Cells.Item(15, 1).Value = "06:52:38"
$ws3.Cells.Item(15, 2).Value = "08:06"
$ws3.Cells.Item(15, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(15, 4).Value = 74
$ws3.Cells.Item(15, 5).Value = "L6203"

$ws3.Cells.Item(18, 1).Value = "06:52:38"
$ws3.Cells.Item(18, 2).Value = "08:34"
$ws3.Cells.Item(18, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(18, 4).Value = 102
$ws3.Cells.Item(18, 5).Value = "L6173"
